$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.804.70'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '2.466.93'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.30%  '
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").Value = '2.467.37'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000177'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '2.916.94'
$ws.Range("D17").Value = '62.804.99'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '2.464.82'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +6.34%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  +17.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '641.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = '0.0₃0984'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.982'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -16.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.26%  '
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("E34").Value = '  -2.10%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  -0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.16%  '
$ws.Range("B42").Value = 'EthereumClassic'
$ws.Range("C42").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '0.0₆0305'
$ws.Range("E45").Value = '  -29.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '153.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.609'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("E51").Value = '  -1.29%  '
